$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting existing C..AI to D..AJ
# (mirrors Excel's native "Insert Column" - it shifts cell content/styles
# right and the new column inherits formatting from the column to its left)
$null = $ws.Columns("C").Insert()

# New column header and value
$ws.Range("C1").Value = "Project Number 2"
$ws.Range("C2").Value = "0000/0304"

# Update existing Project Number value in column B
$ws.Range("B2").Value = "0102/0000"

# Set column C width to match the target layout (closest value this
# host's character->pixel quantization can reach to the authored 16.5703125)
$ws.Columns("C").ColumnWidth = 15.7

# Update the active selection
$null = $ws.Range("B3").Select()
